$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117. This shifts existing rows 117:302 down to 118:303,
# carrying their values and formatting with them.
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new record's data.
$ws.Range("A117").Value2 = 3
$ws.Range("B117").Value2 = "Femacal de La Calera"
$ws.Range("C117").Value2 = "Coquimbo"
$ws.Range("D117").Value2 = 44665
$ws.Range("E117").Value2 = 5
$ws.Range("F117").Value2 = 100112039
$ws.Range("G117").Value2 = "Ciboulette"
$ws.Range("H117").Value2 = "Sin especificar"
$ws.Range("I117").Value2 = "Primera"
$ws.Range("J117").Value2 = 120
$ws.Range("K117").Value2 = 1500
$ws.Range("L117").Value2 = 1500
$ws.Range("M117").Value2 = 1500
$ws.Range("N117").Value2 = '$/docena de atados'
$ws.Range("O117").Value2 = "Provincia de Quillota"
$ws.Range("P117").Value2 = 500
$ws.Range("Q117").Value2 = 3
$ws.Range("R117").Value2 = "Hortaliza"
